$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (12) continuing the "YYYY年" / headcount series,
# matching the formatting (style) already used by the preceding rows.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A12").Value = "2021年"
$ws.Range("B12").Value = 333852
